$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 21.139235
$ws.Range("H2").Value = 63.417705
$ws.Range("I2").Value = 0.1633331201667119
$ws.Range("J2").Value = 0.1633331201667119
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.2296153333333334
$ws.Range("N2").Value = 0.6888460000000001
$ws.Range("O2").Value = 0.04381656765025366
$ws.Range("P2").Value = 0.04381656765025366
$ws.Range("Q2").Value = 4.853892490936667
$ws.Range("R2").Value = 43.68503241843
$ws.Range("S2").Value = 0.007156696709311744
$ws.Range("T2").Value = 0.007156696709311745
# Row 3
$ws.Range("G3").Value = 21.139235
$ws.Range("H3").Value = 63.417705
$ws.Range("I3").Value = 0.1633331201667119
$ws.Range("J3").Value = 0.1633331201667119
$ws.Range("O3").Value = 0.7547076606638542
$ws.Range("P3").Value = 0.7547076606638543
$ws.Range("Q3").Value = 83.60467383454333
$ws.Range("R3").Value = 752.44206451089
$ws.Range("S3").Value = 0.1232687570299473
$ws.Range("T3").Value = 0.1232687570299474
# Row 4
$ws.Range("G4").Value = 21.139235
$ws.Range("H4").Value = 63.417705
$ws.Range("I4").Value = 0.1633331201667119
$ws.Range("J4").Value = 0.1633331201667119
$ws.Range("M4").Value = 1.055809
$ws.Range("N4").Value = 3.167427
$ws.Range("O4").Value = 0.201475771685892
$ws.Range("P4").Value = 0.2014757716858921
$ws.Range("Q4").Value = 22.318994566115
$ws.Range("R4").Value = 200.870951095035
$ws.Range("S4").Value = 0.03290766642745282
$ws.Range("T4").Value = 0.03290766642745283
# Row 5
$ws.Range("I5").Value = 0.5748271090353965
$ws.Range("J5").Value = 0.5748271090353966
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.2296153333333334
$ws.Range("N5").Value = 0.6888460000000001
$ws.Range("O5").Value = 0.04381656765025366
$ws.Range("P5").Value = 0.04381656765025366
$ws.Range("Q5").Value = 17.08256711979711
$ws.Range("R5").Value = 153.743104078174
$ws.Range("S5").Value = 0.02518695091024919
$ws.Range("T5").Value = 0.0251869509102492
# Row 6
$ws.Range("I6").Value = 0.5748271090353965
$ws.Range("J6").Value = 0.5748271090353966
$ws.Range("O6").Value = 0.7547076606638542
$ws.Range("P6").Value = 0.7547076606638543
$ws.Range("S6").Value = 0.4338264227462703
$ws.Range("T6").Value = 0.4338264227462705
# Row 7
$ws.Range("I7").Value = 0.5748271090353965
$ws.Range("J7").Value = 0.5748271090353966
$ws.Range("M7").Value = 1.055809
$ws.Range("N7").Value = 3.167427
$ws.Range("O7").Value = 0.201475771685892
$ws.Range("P7").Value = 0.2014757716858921
$ws.Range("Q7").Value = 78.54844816484034
$ws.Range("R7").Value = 706.9360334835631
$ws.Range("S7").Value = 0.1158137353788769
$ws.Range("T7").Value = 0.1158137353788769
# Row 8
$ws.Range("G8").Value = 4.054539666666667
$ws.Range("H8").Value = 12.163619
$ws.Range("I8").Value = 0.03132755819197652
$ws.Range("J8").Value = 0.03132755819197652
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.2296153333333334
$ws.Range("N8").Value = 0.6888460000000001
$ws.Range("O8").Value = 0.04381656765025366
$ws.Range("P8").Value = 0.04381656765025366
$ws.Range("Q8").Value = 0.9309844770748891
$ws.Range("R8").Value = 8.378860293674002
$ws.Range("S8").Value = 0.001372666072835998
$ws.Range("T8").Value = 0.001372666072835998
# Row 9
$ws.Range("G9").Value = 4.054539666666667
$ws.Range("H9").Value = 12.163619
$ws.Range("I9").Value = 0.03132755819197652
$ws.Range("J9").Value = 0.03132755819197652
$ws.Range("O9").Value = 0.7547076606638542
$ws.Range("P9").Value = 0.7547076606638543
$ws.Range("Q9").Value = 16.03551246678911
$ws.Range("R9").Value = 144.319612201102
$ws.Range("S9").Value = 0.02364314815737736
$ws.Range("T9").Value = 0.02364314815737737
# Row 10
$ws.Range("G10").Value = 4.054539666666667
$ws.Range("H10").Value = 12.163619
$ws.Range("I10").Value = 0.03132755819197652
$ws.Range("J10").Value = 0.03132755819197652
$ws.Range("M10").Value = 1.055809
$ws.Range("N10").Value = 3.167427
$ws.Range("O10").Value = 0.201475771685892
$ws.Range("P10").Value = 0.2014757716858921
$ws.Range("Q10").Value = 4.280819470923666
$ws.Range("R10").Value = 38.527375238313
$ws.Range("S10").Value = 0.006311743961763159
$ws.Range("T10").Value = 0.006311743961763159
# Row 11
$ws.Range("G11").Value = 29.83382566666667
$ws.Range("H11").Value = 89.50147699999999
$ws.Range("I11").Value = 0.2305122126059151
$ws.Range("J11").Value = 0.2305122126059151
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.2296153333333334
$ws.Range("N11").Value = 0.6888460000000001
$ws.Range("O11").Value = 0.04381656765025366
$ws.Range("P11").Value = 0.04381656765025366
$ws.Range("Q11").Value = 6.850303825060223
$ws.Range("R11").Value = 61.652734425542
$ws.Range("S11").Value = 0.01010025395785673
$ws.Range("T11").Value = 0.01010025395785673
# Row 12
$ws.Range("G12").Value = 29.83382566666667
$ws.Range("H12").Value = 89.50147699999999
$ws.Range("I12").Value = 0.2305122126059151
$ws.Range("J12").Value = 0.2305122126059151
$ws.Range("O12").Value = 0.7547076606638542
$ws.Range("P12").Value = 0.7547076606638543
$ws.Range("Q12").Value = 117.9913683772518
$ws.Range("R12").Value = 1061.922315395266
$ws.Range("S12").Value = 0.1739693327302592
$ws.Range("T12").Value = 0.1739693327302592
# Row 13
$ws.Range("G13").Value = 29.83382566666667
$ws.Range("H13").Value = 89.50147699999999
$ws.Range("I13").Value = 0.2305122126059151
$ws.Range("J13").Value = 0.2305122126059151
$ws.Range("M13").Value = 1.055809
$ws.Range("N13").Value = 3.167427
$ws.Range("O13").Value = 0.201475771685892
$ws.Range("P13").Value = 0.2014757716858921
$ws.Range("Q13").Value = 31.49882164329767
$ws.Range("R13").Value = 283.489394789679
$ws.Range("S13").Value = 0.04644262591779915
$ws.Range("T13").Value = 0.04644262591779916
